# Apply the asset.xlsx edit described by the commit "coding the fast tool function"
#
# Summary of the change:
#  - Comment on C2 gains extra lines ("在线:1" / "离线：2") and loses the bold
#    run that used to prefix it.
#  - A new comment is added on U2 ("作者:" + newline + "用户ID").
#  - sharedStrings gains three new strings: 192.168.10.4, 192.168.10.5 and
#    所有者ID, while the old "192.168.10.7" string is repointed to
#    "192.168.10.3".
#  - Column U (header value 20) is appended to the header row.
#  - Row 3 gains a B3 value (22) and a U3 value (1); F3 changes from 6 to 1.
#  - Two brand-new data rows (4 and 5) are added, mirroring row 3's layout,
#    for hosts 192.168.10.4 and 192.168.10.5.
#  - The old, mostly-empty row 8 (B8:D8 styled blanks) is removed.
#  - The sheet view now starts scrolled to column G with U2 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Comment text update on C2 (loses the leading bold run, gains two lines)
# ---------------------------------------------------------------------------
$c2Comment = $ws.Range("C2").Comment
$c2Comment.Text("必填项`n在线:1`n离线：2")

# ---------------------------------------------------------------------------
# 2. Header row (row 1) gains column U
# ---------------------------------------------------------------------------
$ws.Range("U1").Value = 20

# ---------------------------------------------------------------------------
# 3. Row 3 edits: re-point A3's IP, add B3, change F3, add U3
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "192.168.10.3"
$ws.Range("B3").Value = 22
$ws.Range("F3").Value = 1
$ws.Range("U3").Value = 1

# ---------------------------------------------------------------------------
# 4. New rows 4 and 5, cloned from row 3's cell formatting
# ---------------------------------------------------------------------------
foreach ($pair in @(@(4, "192.168.10.4"), @(5, "192.168.10.5"))) {
    $r = $pair[0]
    $ip = $pair[1]

    $ws.Range("A3").Copy($ws.Range("A$r"))
    $ws.Range("C3").Copy($ws.Range("C$r"))
    $ws.Range("D3").Copy($ws.Range("D$r"))
    $ws.Range("E3").Copy($ws.Range("E$r"))
    $ws.Range("F3").Copy($ws.Range("F$r"))
    $ws.Range("G3").Copy($ws.Range("G$r"))
    $ws.Range("S3").Copy($ws.Range("S$r"))

    $ws.Range("A$r").Value = $ip
    $ws.Range("B$r").Value = 22
    $ws.Range("F$r").Value = 1
    $ws.Range("U$r").Value = 1
}

# ---------------------------------------------------------------------------
# 5. New cell + comment on U2 ("作者:" bold run, then a plain run with the
#    newline + "用户ID" - the engine stores comments as plain text, so we
#    just concatenate the runs' text). Style matches the rest of the header
#    row (same as T2).
# ---------------------------------------------------------------------------
$ws.Range("T2").Copy($ws.Range("U2"))
$ws.Range("U2").Value = "所有者ID"
$ws.Range("U2").AddComment("作者:`n用户ID")

# ---------------------------------------------------------------------------
# 6. Drop the old, mostly-empty row 8
# ---------------------------------------------------------------------------
$ws.Rows.Item(8).Delete()

# ---------------------------------------------------------------------------
# 7. Sheet view: scroll to column G, select U2
# ---------------------------------------------------------------------------
$ws.Range("U2").Select()
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
